# Apply the vocabulary.xlsx update:
#  - refresh the "dct:modified" timestamp (B21)
#  - simplify the skos:broader column mapping spec (F23)
#  - rename the placeholder demo term (B24)
#  - append the newly generated vocabulary terms (rows 25-37)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updates to existing cells -------------------------------------------------
$ws.Range("B21").Value = "2023-08-17T11:19:58+00:00"
$ws.Range("F23").Value = 'skos:broader(separator=",")'
$ws.Range("B24").Value = "grant type"

# --- new term rows appended at the bottom of the sheet -------------------------
# Each entry is: Identifier, skos:prefLabel, skos:broader (comma separated ids)
$newTerms = @(
    @("vocab:1001", "new",            "vocab.1001"),
    @("vocab:1002", "repurposed",     "vocab.1001"),
    @("vocab:1003", "supplemented",   "vocab.1001"),
    @("vocab:1004", "unspecified",    "vocab.1001,vocab.1007"),
    @("vocab:1005", "not applicable", "vocab.1001,vocab.1007"),
    @("vocab:1006", "not known",      "vocab.1001,vocab.1007"),
    @("vocab:1007", "age group",      ""),
    @("vocab:1008", "adolescent",     "vocab.1007"),
    @("vocab:1009", "adult",          "vocab.1007"),
    @("vocab:1010", "fail elderly",   "vocab.1007"),
    @("vocab:1011", "child",          "vocab.1007"),
    @("vocab:1012", "infant",         "vocab.1007"),
    @("vocab:1013", "newborn",        "vocab.1007")
)

$row = 25
foreach ($term in $newTerms) {
    $ws.Cells.Item($row, 1).Value = $term[0]
    $ws.Cells.Item($row, 2).Value = $term[1]
    if ($term[2] -ne "") {
        $ws.Cells.Item($row, 3).Value = $term[2]
    }
    $row++
}
